$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 44.1
$ws.Range("I28").Value = 44.1
$ws.Range("K28").Value = 44.1
$ws.Range("M28").Value = 440.9
$ws.Range("H32").Value = 8553.944
$ws.Range("I32").Value = 7230.5
$ws.Range("K32").Value = 7230.5
$ws.Range("M32").Value = -6904.5
$ws.Range("H69").Value = 12000
$ws.Range("J69").Value = 12000
$ws.Range("L69").Value = 36000
$ws.Range("N69").Value = -37748
$ws.Range("H72").Value = 12000
$ws.Range("J72").Value = 12000
$ws.Range("L72").Value = 108000
$ws.Range("N72").Value = -116736
$ws.Range("H76").Value = 5000
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 5000
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H97").Value = 8180
$ws.Range("J97").Value = 8180
$ws.Range("L97").Value = 24540
$ws.Range("N97").Value = -25532
$ws.Range("I137").Value = 4902
$ws.Range("J137").Value = 4850.3335
$ws.Range("K137").Value = 14706
$ws.Range("L137").Value = 14551.0005
$ws.Range("M137").Value = -12156
$ws.Range("N137").Value = -19651.0005
$ws.Range("H138").Value = 7388.684
$ws.Range("J138").Value = 7521.3887
$ws.Range("L138").Value = 22564.1661
$ws.Range("N138").Value = -32844.1661
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13858.714
$ws.Range("I32").Value = 12183.818
$ws.Range("K32").Value = 12183.818
$ws.Range("M32").Value = -11896.818
$ws.Range("H52").Value = 49990
$ws.Range("J52").Value = 49990
$ws.Range("L52").Value = 49990
$ws.Range("N52").Value = -50626
$ws.Range("H61").Value = 4469.7144
$ws.Range("I61").Value = 4257.6
$ws.Range("K61").Value = 4257.6
$ws.Range("M61").Value = -4045.6
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("H110").Value = 1379.2
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("H136").Value = 4469.7144
$ws.Range("I136").Value = 4257.6
$ws.Range("K136").Value = 12772.8
$ws.Range("M136").Value = -10222.8
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()
$ws.Range("N110").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6163.5454
$ws.Range("I86").Value = 2500
$ws.Range("J86").Value = 6977.6665
$ws.Range("K86").Value = 2500
$ws.Range("L86").Value = 6977.6665
$ws.Range("M86").Value = -1377
$ws.Range("N86").Value = -9223.666499999999
$ws.Range("H89").Value = 6163.5454
$ws.Range("I89").Value = 2500
$ws.Range("J89").Value = 6977.6665
$ws.Range("K89").Value = 12500
$ws.Range("L89").Value = 34888.3325
$ws.Range("M89").Value = -6884
$ws.Range("N89").Value = -46120.3325
$ws.Range("H105").Value = 2271.4285
$ws.Range("I105").Value = 2271.4285
$ws.Range("K105").Value = 2271.4285
$ws.Range("M105").Value = -524.4285

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 100
$ws.Range("M7").Value = 13
$ws.Range("H31").Value = 43124
$ws.Range("I31").Value = 23333
$ws.Range("K31").Value = 23333
$ws.Range("M31").Value = -23038
$ws.Range("H34").Value = 43124
$ws.Range("I34").Value = 23333
$ws.Range("K34").Value = 23333
$ws.Range("M34").Value = -23131

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.42857
$ws.Range("I2").Value = 56.42857
$ws.Range("J2").Value = 44.42857
$ws.Range("K2").Value = 338.57142
$ws.Range("L2").Value = 266.57142
$ws.Range("M2").Value = -225.57142
$ws.Range("N2").Value = -492.57142
$ws.Range("H4").Value = 166871.75
$ws.Range("I4").Value = 325
$ws.Range("J4").Value = 285833.72
$ws.Range("K4").Value = 975
$ws.Range("L4").Value = 857501.1599999999
$ws.Range("M4").Value = -863
$ws.Range("N4").Value = -857725.1599999999
$ws.Range("H131").Value = 4209.3076
$ws.Range("J131").Value = 3394.1
$ws.Range("L131").Value = 10182.3
$ws.Range("N131").Value = -20262.3
$ws.Range("H134").Value = 957.6667
$ws.Range("I134").Value = 957.6667
$ws.Range("K134").Value = 2873.0001
$ws.Range("M134").Value = 2196.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 850
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405
$ws.Range("H27").Value = 850
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593
$ws.Range("H46").Value = 975
$ws.Range("I46").Value = 800
$ws.Range("J46").Value = 1033.3334
$ws.Range("K46").Value = 800
$ws.Range("L46").Value = 1033.3334
$ws.Range("M46").Value = -612
$ws.Range("N46").Value = -1409.3334
$ws.Range("H54").Value = 45001
$ws.Range("J54").Value = 45001
$ws.Range("L54").Value = 45001
$ws.Range("N54").Value = -46289
$ws.Range("H82").Value = 2054.0908
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 2054.0908
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H132").Value = 70000
$ws.Range("I132").Value = 80000
$ws.Range("K132").Value = 240000
$ws.Range("M132").Value = -237470
$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 20000000
$ws.Range("J46").Value = 20000000
$ws.Range("L46").Value = 20000000
$ws.Range("N46").Value = -20000462
$ws.Range("H81").Value = 425.83334
$ws.Range("I81").Value = 418.75
$ws.Range("J81").Value = 440
$ws.Range("K81").Value = 837.5
$ws.Range("L81").Value = 880
$ws.Range("M81").Value = 223.5
$ws.Range("N81").Value = -3002
$ws.Range("H84").Value = 425.83334
$ws.Range("I84").Value = 418.75
$ws.Range("J84").Value = 440
$ws.Range("K84").Value = 4187.5
$ws.Range("L84").Value = 4400
$ws.Range("M84").Value = 1116.5
$ws.Range("N84").Value = -15008
$ws.Range("H132").Value = 5441.478
$ws.Range("I132").Value = 4730.8335
$ws.Range("J132").Value = 7999.8
$ws.Range("K132").Value = 14192.5005
$ws.Range("L132").Value = 23999.4
$ws.Range("M132").Value = -11662.5005
$ws.Range("N132").Value = -29059.4
$ws.Range("H134").Value = 20000000
$ws.Range("J134").Value = 20000000
$ws.Range("L134").Value = 60000000
$ws.Range("N134").Value = -60005070
